# Season-record columns: Wins, Losses, Ties
# Adds three new columns (AD, AE, AF) to the player stats sheet, holding
# the team's season record, repeated for every player row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting (bold, centered, bordered) from an existing
# header cell so the new headers match the rest of row 1's style.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows (2 through 52): same season record for every player ---
$wins = 79
$losses = 83
$ties = 0

for ($r = 2; $r -le 52; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins
    $ws.Cells.Item($r, 31).Value = $losses
    $ws.Cells.Item($r, 32).Value = $ties
}
